# Add a new "Link" message type.
# - Inserts a new "Link" column (H) into the User_Initiated_Messages and
#   Follow_Up_Messages header tables (pushing the existing "Follow Ups"
#   column from H to I).
# - Adds a new "link" message row to User_Initiated_Messages with a real
#   hyperlink to the choicenotchance.org.nz site.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# User_Initiated_Messages
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("User_Initiated_Messages")

# Make room for the new "Link" column before the existing "Follow Ups"
# column (old H -> new I), inheriting formatting from the left neighbour
# the same way Excel's own "Insert Sheet Columns" does.
$ws3.Columns("H:H").Insert()

$ws3.Range("H1").Value = "Link"

# New row describing the "link" message type.
$ws3.Range("A6").Value = 4
$ws3.Range("B6").Value = "link"
$ws3.Range("C6").Value = "Message;Link"
$ws3.Range("D6").Value = "Try these links for help"

# The link itself, as a real hyperlink (Excel auto-applies the built-in
# Hyperlink style and uses the target URL as the visible text).
$ws3.Hyperlinks.Add($ws3.Range("H6"), "https://www.choicenotchance.org.nz/") | Out-Null

$ws3.Range("F12").Select() | Out-Null

# ---------------------------------------------------------------------
# Follow_Up_Messages
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Follow_Up_Messages")

$ws4.Columns("H:H").Insert()

$ws4.Range("H1").Value = "Link"

$ws4.Activate() | Out-Null
$ws4.Range("H5").Select() | Out-Null
